# Auto-generated edit script: updates market-price derived columns (H-N)
# across leve-profit sheets, per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 33946.5
$ws.Range("I20").Value = 5875
$ws.Range("J20").Value = 62018
$ws.Range("K20").Value = 5875
$ws.Range("L20").Value = 62018
$ws.Range("M20").Value = -5645
$ws.Range("N20").Value = -62478
$ws.Range("H35").Value = 33946.5
$ws.Range("I35").Value = 5875
$ws.Range("J35").Value = 62018
$ws.Range("K35").Value = 5875
$ws.Range("L35").Value = 62018
$ws.Range("M35").Value = -5496
$ws.Range("N35").Value = -62776
$ws.Range("H112").Value = 1266.5625
$ws.Range("I112").Value = 750
$ws.Range("K112").Value = 2250
$ws.Range("M112").Value = -1142
$ws.Range("H121").Value = 1477.5
$ws.Range("J121").Value = 1477.5
$ws.Range("L121").Value = 4432.5
$ws.Range("N121").Value = -7926.5
$ws.Range("H131").Value = 1096.5
$ws.Range("I131").Value = 720.625
$ws.Range("K131").Value = 2161.875
$ws.Range("M131").Value = 2878.125
$ws.Range("H132").Value = 4465400.5
$ws.Range("I132").Value = 5715372.5
$ws.Range("J132").Value = 1215.1428
$ws.Range("K132").Value = 17146117.5
$ws.Range("L132").Value = 3645.4284
$ws.Range("M132").Value = -17143587.5
$ws.Range("N132").Value = -8705.428400000001
$ws.Range("H137").Value = 1401.8718
$ws.Range("I137").Value = 1264.2903
$ws.Range("J137").Value = 1935
$ws.Range("K137").Value = 3792.8709
$ws.Range("L137").Value = 5805
$ws.Range("M137").Value = -1242.8709
$ws.Range("N137").Value = -10905
$ws.Range("H141").Value = 4546.037
$ws.Range("I141").Value = 1434.375
$ws.Range("K141").Value = 4303.125
$ws.Range("M141").Value = 876.875

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1041.8529
$ws.Range("I2").Value = 776.6667
$ws.Range("J2").Value = 1470.2307
$ws.Range("K2").Value = 776.6667
$ws.Range("L2").Value = 1470.2307
$ws.Range("M2").Value = -663.6667
$ws.Range("N2").Value = -1696.2307
$ws.Range("H45").Value = 1212.1765
$ws.Range("I45").Value = 1220.6364
$ws.Range("K45").Value = 1220.6364
$ws.Range("M45").Value = -843.6364000000001
$ws.Range("H57").Value = 8000
$ws.Range("I57").Value = 8000
$ws.Range("K57").Value = 8000
$ws.Range("M57").Value = -7516
$ws.Range("H61").Value = 1560.8518
$ws.Range("I61").Value = 1255.9584
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1255.9584
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -1043.9584
$ws.Range("N61").Value = -4424
$ws.Range("H116").Value = 1041.8529
$ws.Range("I116").Value = 776.6667
$ws.Range("J116").Value = 1470.2307
$ws.Range("K116").Value = 776.6667
$ws.Range("L116").Value = 1470.2307
$ws.Range("M116").Value = 1517.3333
$ws.Range("N116").Value = -6058.2307
$ws.Range("H136").Value = 1560.8518
$ws.Range("I136").Value = 1255.9584
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 3767.8752
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -1217.8752
$ws.Range("N136").Value = -17100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1041.8529
$ws.Range("I3").Value = 776.6667
$ws.Range("J3").Value = 1470.2307
$ws.Range("K3").Value = 776.6667
$ws.Range("L3").Value = 1470.2307
$ws.Range("M3").Value = -662.6667
$ws.Range("N3").Value = -1698.2307
$ws.Range("H20").Value = 1953.6666
$ws.Range("I20").Value = 2088.0625
$ws.Range("J20").Value = 1684.875
$ws.Range("K20").Value = 2088.0625
$ws.Range("L20").Value = 1684.875
$ws.Range("M20").Value = -1841.0625
$ws.Range("N20").Value = -2178.875
$ws.Range("H44").Value = 15333.333
$ws.Range("J44").Value = 15333.333
$ws.Range("L44").Value = 15333.333
$ws.Range("N44").Value = -16327.333
$ws.Range("H94").Value = 1149.9375
$ws.Range("I94").Value = 907.0714
$ws.Range("K94").Value = 907.0714
$ws.Range("M94").Value = -456.0714
$ws.Range("H105").Value = 2155
$ws.Range("I105").Value = 2086
$ws.Range("K105").Value = 2086
$ws.Range("M105").Value = -339

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5265551.5
$ws.Range("I31").Value = 2453.3333
$ws.Range("J31").Value = 40002000
$ws.Range("K31").Value = 2453.3333
$ws.Range("L31").Value = 40002000
$ws.Range("M31").Value = -2158.3333
$ws.Range("N31").Value = -40002590
$ws.Range("H34").Value = 5265551.5
$ws.Range("I34").Value = 2453.3333
$ws.Range("J34").Value = 40002000
$ws.Range("K34").Value = 2453.3333
$ws.Range("L34").Value = 40002000
$ws.Range("M34").Value = -2251.3333
$ws.Range("N34").Value = -40002404
$ws.Range("H54").Value = 16915
$ws.Range("J54").Value = 16915
$ws.Range("L54").Value = 16915
$ws.Range("N54").Value = -18231
$ws.Range("H58").Value = 1129.8695
$ws.Range("I58").Value = 965.8
$ws.Range("J58").Value = 1437.5
$ws.Range("K58").Value = 965.8
$ws.Range("L58").Value = 1437.5
$ws.Range("M58").Value = -762.8
$ws.Range("N58").Value = -1843.5
$ws.Range("H99").Value = 2245
$ws.Range("I99").Value = 1750
$ws.Range("J99").Value = 2575
$ws.Range("K99").Value = 1750
$ws.Range("L99").Value = 2575
$ws.Range("M99").Value = -252
$ws.Range("N99").Value = -5571
$ws.Range("H126").Value = 2245
$ws.Range("I126").Value = 1750
$ws.Range("J126").Value = 2575
$ws.Range("K126").Value = 5250
$ws.Range("L126").Value = 7725
$ws.Range("M126").Value = -2780
$ws.Range("N126").Value = -12665
$ws.Range("H132").Value = 3458.647
$ws.Range("I132").Value = 2914.2144
$ws.Range("J132").Value = 5999.3335
$ws.Range("K132").Value = 8742.643199999999
$ws.Range("L132").Value = 17998.0005
$ws.Range("M132").Value = -6212.643199999999
$ws.Range("N132").Value = -23058.0005
$ws.Range("H134").Value = 813.15216
$ws.Range("I134").Value = 768.0732
$ws.Range("K134").Value = 2304.2196
$ws.Range("M134").Value = 230.7803999999996
$ws.Range("H136").Value = 1129.8695
$ws.Range("I136").Value = 965.8
$ws.Range("J136").Value = 1437.5
$ws.Range("K136").Value = 2897.4
$ws.Range("L136").Value = 4312.5
$ws.Range("M136").Value = -347.3999999999996
$ws.Range("N136").Value = -9412.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 14680.167
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 14680.167
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 44040.501
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -45392.501
$ws.Range("H139").Value = 2149.1875
$ws.Range("I139").Value = 709.6667
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 2129.0001
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = 3010.9999
$ws.Range("N139").Value = -22280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18894762
$ws.Range("I70").Value = 25504220
$ws.Range("J70").Value = 10600
$ws.Range("K70").Value = 25504220
$ws.Range("L70").Value = 10600
$ws.Range("M70").Value = -25503950
$ws.Range("N70").Value = -11140
$ws.Range("H73").Value = 18894762
$ws.Range("I73").Value = 25504220
$ws.Range("J73").Value = 10600
$ws.Range("K73").Value = 25504220
$ws.Range("L73").Value = 10600
$ws.Range("M73").Value = -25503284
$ws.Range("N73").Value = -12472
$ws.Range("H99").Value = 1510.5
$ws.Range("I99").Value = 1510.5
$ws.Range("K99").Value = 1510.5
$ws.Range("M99").Value = 735.5
$ws.Range("H126").Value = 5557161.5
$ws.Range("I126").Value = 1678.6666
$ws.Range("J126").Value = 11112644
$ws.Range("K126").Value = 5035.9998
$ws.Range("L126").Value = 33337932
$ws.Range("M126").Value = -2565.9998
$ws.Range("N126").Value = -33342872

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1862.35
$ws.Range("I40").Value = 1820.4117
$ws.Range("K40").Value = 1820.4117
$ws.Range("M40").Value = -1684.4117
$ws.Range("H132").Value = 3242.4375
$ws.Range("I132").Value = 3250.5789
$ws.Range("K132").Value = 9751.736699999999
$ws.Range("M132").Value = -7221.736699999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1588.6923
$ws.Range("I122").Value = 1637.75
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4913.25
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2463.25
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 33669.332
$ws.Range("I126").Value = 33669.332
$ws.Range("K126").Value = 101007.996
$ws.Range("M126").Value = -98537.99600000001
$ws.Range("H132").Value = 1788
$ws.Range("I132").Value = 1241.9286
$ws.Range("K132").Value = 3725.7858
$ws.Range("M132").Value = -1195.7858
$ws.Range("H136").Value = 7043.3
$ws.Range("I136").Value = 7043.3
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 21129.9
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -18579.9
$ws.Range("N136").ClearContents()
